$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 5.432019178035748
$ws.Range("D2").Value = 14.37994989851497
$ws.Range("E2").Value = 10.96708249694544
$ws.Range("F2").Value = 70.8277802207096
$ws.Range("G2").Value = 3.77315679165562
$ws.Range("I2").Value = 57.37688974334848
$ws.Range("L2").Value = 10.30711645587346
$ws.Range("N2").Value = 18.99769123883691
$ws.Range("C3").Value = 5.33238946464488
$ws.Range("D3").Value = 14.20838952785987
$ws.Range("E3").Value = 11.00335940245331
$ws.Range("F3").Value = 69.33793411521458
$ws.Range("G3").Value = 3.78628358837993
$ws.Range("I3").Value = 56.39056763191363
$ws.Range("L3").Value = 10.32880045711248
$ws.Range("N3").Value = 18.40031508502705
$ws.Range("C4").Value = 5.269859457617176
$ws.Range("D4").Value = 14.10867205384778
$ws.Range("E4").Value = 11.02831929866602
$ws.Range("F4").Value = 68.44323962109851
$ws.Range("G4").Value = 3.794683062915857
$ws.Range("I4").Value = 55.80438599385968
$ws.Range("L4").Value = 10.34498935555675
$ws.Range("N4").Value = 18.02485520896364
$ws.Range("C5").Value = 5.244044124759378
$ws.Range("D5").Value = 14.06946847375139
$ws.Range("E5").Value = 11.03916136935733
$ws.Range("F5").Value = 68.08410711132531
$ws.Range("G5").Value = 3.798192345767762
$ws.Range("I5").Value = 55.57061343284588
$ws.Range("L5").Value = 10.35230219779425
$ws.Range("N5").Value = 17.86990355188116
$ws.Range("C6").Value = 5.239737570217142
$ws.Range("D6").Value = 14.06304571961459
$ws.Range("E6").Value = 11.04100206297787
$ws.Range("F6").Value = 68.02481530256783
$ws.Range("G6").Value = 3.798780309588909
$ws.Range("I6").Value = 55.53210987006684
$ws.Range("L6").Value = 10.35355950979229
$ws.Range("N6").Value = 17.84406337566596
$ws.Range("C7").Value = 5.269512644643617
$ws.Range("D7").Value = 14.10813752009782
$ws.Range("E7").Value = 11.0284628089587
$ws.Range("F7").Value = 68.43837357469987
$ws.Range("G7").Value = 3.794730039022434
$ws.Range("I7").Value = 55.8012123231704
$ws.Range("L7").Value = 10.34508509053109
$ws.Range("N7").Value = 18.02277304766463
$ws.Range("C8").Value = 5.397952731411674
$ws.Range("D8").Value = 14.31964132930881
$ws.Range("E8").Value = 10.97903052001517
$ws.Range("F8").Value = 70.31016472982139
$ws.Range("G8").Value = 3.777613069900291
$ws.Range("I8").Value = 57.03290791136932
$ws.Range("L8").Value = 10.31399209573145
$ws.Range("N8").Value = 18.79364780656866
$ws.Range("C9").Value = 5.63876833109755
$ws.Range("D9").Value = 14.77804291066461
$ws.Range("E9").Value = 10.90361980359458
$ws.Range("F9").Value = 74.12171597412453
$ws.Range("G9").Value = 3.746692455716552
$ws.Range("I9").Value = 59.59241257954898
$ws.Range("L9").Value = 10.27615572341387
$ws.Range("N9").Value = 20.2273683202997
$ws.Range("C10").Value = 5.808627047943533
$ws.Range("D10").Value = 15.13999543410459
$ws.Range("E10").Value = 10.86164768903911
$ws.Range("F10").Value = 76.98386661283966
$ws.Range("G10").Value = 3.725518037624417
$ws.Range("I10").Value = 61.54714058716828
$ws.Range("L10").Value = 10.26291940535079
$ws.Range("N10").Value = 21.22223697909767
$ws.Range("C11").Value = 5.88431100790652
$ws.Range("D11").Value = 15.30977916409777
$ws.Range("E11").Value = 10.84554068159729
$ws.Range("F11").Value = 78.29449542027703
$ws.Range("G11").Value = 3.716204406039463
$ws.Range("I11").Value = 62.44964211473809
$ws.Range("L11").Value = 10.26016049763678
$ws.Range("N11").Value = 21.66018057919901
$ws.Range("C12").Value = 5.91273882299694
$ws.Range("D12").Value = 15.37477905394122
$ws.Range("E12").Value = 10.83987690490871
$ws.Range("F12").Value = 78.79165728588413
$ws.Range("G12").Value = 3.712722064660348
$ws.Range("I12").Value = 62.79307245729171
$ws.Range("L12").Value = 10.25959328643976
$ws.Range("N12").Value = 21.82377585681918
$ws.Range("C13").Value = 5.906626760583528
$ws.Range("D13").Value = 15.36074924082907
$ws.Range("E13").Value = 10.84107722637058
$ws.Range("F13").Value = 78.68455239002708
$ws.Range("G13").Value = 3.713470090079163
$ws.Range("I13").Value = 62.71903787351015
$ws.Range("L13").Value = 10.2596940717787
$ws.Range("N13").Value = 21.78864458690801
$ws.Range("C14").Value = 5.886654452309262
$ws.Range("D14").Value = 15.31511276413469
$ws.Range("E14").Value = 10.84506595304095
$ws.Range("F14").Value = 78.33538219626183
$ws.Range("G14").Value = 3.715917026950037
$ws.Range("I14").Value = 62.47786413950131
$ws.Range("L14").Value = 10.2601042137446
$ws.Range("N14").Value = 21.67368539489599
$ws.Range("C15").Value = 5.874390510357419
$ws.Range("D15").Value = 15.28725018157197
$ws.Range("E15").Value = 10.8475660782476
$ws.Range("F15").Value = 78.12160577313388
$ws.Range("G15").Value = 3.717421605816716
$ws.Range("I15").Value = 62.33034920267111
$ws.Range("L15").Value = 10.26041787791934
$ws.Range("N15").Value = 21.60297336126124
$ws.Range("C16").Value = 5.803648651160008
$ws.Range("D16").Value = 15.12900023081763
$ws.Range("E16").Value = 10.86276093394254
$ws.Range("F16").Value = 76.89835396568864
$ws.Range("G16").Value = 3.726133008653255
$ws.Range("I16").Value = 61.488406773121
$ws.Range("L16").Value = 10.26316606818449
$ws.Range("N16").Value = 21.19330956972086
$ws.Range("C17").Value = 5.759841550375415
$ws.Range("D17").Value = 15.03321076014114
$ws.Range("E17").Value = 10.87285166578735
$ws.Range("F17").Value = 76.14985520502327
$ws.Range("G17").Value = 3.731557893446692
$ws.Range("I17").Value = 60.97513290380904
$ws.Range("L17").Value = 10.26569349160005
$ws.Range("N17").Value = 20.93814219015166
$ws.Range("C18").Value = 5.734495808646508
$ws.Range("D18").Value = 14.97859992446942
$ws.Range("E18").Value = 10.87893629091508
$ws.Range("F18").Value = 75.72017159862972
$ws.Range("G18").Value = 3.734708223949952
$ws.Range("I18").Value = 60.68117561248639
$ws.Range("L18").Value = 10.26745375574026
$ws.Range("N18").Value = 20.79000725568141
$ws.Range("C19").Value = 5.725888712569485
$ws.Range("D19").Value = 14.96019383336798
$ws.Range("E19").Value = 10.88104446745451
$ws.Range("F19").Value = 75.57484333479348
$ws.Range("G19").Value = 3.73578007685061
$ws.Range("I19").Value = 60.58187127160402
$ws.Range("L19").Value = 10.26810214893005
$ws.Range("N19").Value = 20.73962067985785
$ws.Range("C20").Value = 5.764520349768115
$ws.Range("D20").Value = 15.04335776994106
$ws.Range("E20").Value = 10.87174839455799
$ws.Range("F20").Value = 76.22945082397973
$ws.Range("G20").Value = 3.730977301272129
$ws.Range("I20").Value = 61.02964269254235
$ws.Range("L20").Value = 10.26539265522109
$ws.Range("N20").Value = 20.96544799483449
$ws.Range("C21").Value = 5.892527129407968
$ws.Range("D21").Value = 15.32849837882466
$ws.Range("E21").Value = 10.8438824949486
$ws.Range("F21").Value = 78.43792167114036
$ws.Range("G21").Value = 3.715197104429218
$ws.Range("I21").Value = 62.5486591611259
$ws.Range("L21").Value = 10.25997071735805
$ws.Range("N21").Value = 21.70751365554064
$ws.Range("C22").Value = 5.974833081525459
$ws.Range("D22").Value = 15.51895775614673
$ws.Range("E22").Value = 10.82821255722245
$ws.Range("F22").Value = 79.88611074695183
$ws.Range("G22").Value = 3.705142737028785
$ws.Range("I22").Value = 63.55107119365304
$ws.Range("L22").Value = 10.25921458415219
$ws.Range("N22").Value = 22.17935961385674
$ws.Range("C23").Value = 5.931029840813487
$ws.Range("D23").Value = 15.41694081853078
$ws.Range("E23").Value = 10.83634119766232
$ws.Range("F23").Value = 79.1128634935535
$ws.Range("G23").Value = 3.710485697192249
$ws.Range("I23").Value = 63.01525805303074
$ws.Range("L23").Value = 10.25936027710262
$ws.Range("N23").Value = 21.92877110911181
$ws.Range("C24").Value = 5.762405564029657
$ws.Range("D24").Value = 15.03876887224643
$ws.Range("E24").Value = 10.87224630143798
$ws.Range("F24").Value = 76.19346362697257
$ws.Range("G24").Value = 3.731239688806049
$ws.Range("I24").Value = 61.00499527701161
$ws.Range("L24").Value = 10.26552770707653
$ws.Range("N24").Value = 20.95310750188673
$ws.Range("C25").Value = 5.574841451092836
$ws.Range("D25").Value = 14.64950132746626
$ws.Range("E25").Value = 10.92168464148548
$ws.Range("F25").Value = 73.07802671780179
$ws.Range("G25").Value = 3.754781139305293
$ws.Range("I25").Value = 58.88596542044987
$ws.Range("L25").Value = 10.28386766157655
$ws.Range("N25").Value = 19.84905939529497
